# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# to match the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.847.19"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "2.823.90"
$ws.Range("E3").Value = "  +1.68%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.58%  "

$ws.Range("E7").Value = "  +4.73%  "

$ws.Range("E9").Value = "  +4.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.68%  "

$ws.Range("E12").Value = "  +1.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.92%  "

$ws.Range("E14").Value = "  +1.22%  "

$ws.Range("D15").Value = "3.265.83"
$ws.Range("E15").Value = "  +1.87%  "

$ws.Range("D16").Value = "2.824.77"
$ws.Range("E16").Value = "  +2.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.887"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").Value = "51.746.97"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.18%  "

$ws.Range("D22").Value = "0.0₃0989"
$ws.Range("E22").Value = "  +1.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.76%  "

$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.52%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "

$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("E30").Value = "  -1.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "50.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.15%  "

$ws.Range("E32").Value = "  -3.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0455"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +26.90%  "

$ws.Range("E34").Value = "  +4.47%  "

$ws.Range("E35").Value = "  +6.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0826"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.37%  "

$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("E38").Value = "  -1.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.79"
$ws.Range("D41").Style = "Normal"

$ws.Range("E42").Value = "  +2.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "125.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.21%  "

$ws.Range("E44").Value = "  +1.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.50%  "

$ws.Range("D46").Value = "2.087.68"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.934"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "

